# Bergey model completed. Ready for QED adaption.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Update column C (Vdc) values, rows 2-18 ---
$cVals = @{
    2  = 326
    3  = 422
    4  = 520
    5  = 558
    6  = 560
    7  = 533
    8  = 510
    9  = 490
    10 = 475
    11 = 458
    12 = 453
    13 = 459
    14 = 480
    15 = 507
    16 = 533
    17 = 554
    18 = 563
}

foreach ($row in $cVals.Keys) {
    $ws.Cells.Item($row, 3).Value = $cVals[$row]
}

# --- D4 (Power) carries a tiny floating point recalculation difference ---
# Writing a new number straight into D4 causes the host to drop the cell's
# "quotePrefix" style flag (it reassigns a brand new style), which would
# needlessly perturb the stylesheet. Avoid that by resetting D4 to the
# default format first (so there is no quote-prefixed style to disturb),
# writing the new value, and then restoring the original formatting by
# copying it over from an untouched sibling cell that still carries it.
$ws.Cells.Item(4, 4).ClearFormats()
$ws.Cells.Item(4, 4).Value = 3798.9507282711802
$ws.Cells.Item(2, 4).Copy()
$ws.Cells.Item(4, 4).PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update active selection on the worksheet ---
$ws.Range("C3").Select()
